# Rename a handful of "Class Name" entries in the Classes lookup sheet
# (shortened/clarified labels) and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = "Turtle Ships & Thirisadai"
$ws.Range("B15").Value = "Stone Defense & Harb."
$ws.Range("B19").Value = "High Pierce Arm. Siege"
$ws.Range("B10").Value = "Mounted Units (excl. Camels)"
$ws.Range("B18").Value = "Ships (excl. fishing ships)"

# Update the selected cell to match the saved view state.
$ws.Range("B11").Select()
